$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns (I and J), using the same bold/
# bordered/centered header style as the existing headers (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for columns I ("I0") and J ("IF") for rows 2-39.
$IValues = @(6,3,7,8,7,4,8,8,5,8,5,8,9,9,9,7,6,8,7,7,6,7,9,4,2,7,4,5,9,8,9,8,9,9,8,8,7,8)
$JValues = @(7,5,7,9,7,5,9,9,5,9,5,8,9,9,9,7,7,8,8,9,6,7,9,5,4,8,5,6,9,8,9,9,9,9,8,9,7,8)

for ($i = 0; $i -lt $IValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $IValues[$i]
    $ws.Cells.Item($row, 10).Value = $JValues[$i]
}
